$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the shared text in column A (rows 2-32) to change "Slashing" to "Run Slashing"
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = "0_Fallen_Angels_Run Slashing_"
}

# Re-enter the formula in C13 explicitly, which causes Excel to break it out of
# the shared formula group (matches the diff where C13's <f> becomes non-shared)
$ws.Range("C13").Formula = "=A13&B13"

# Update selection to reflect C13 being active cell (matches diff)
$ws.Range("C13").Select()
